$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Update existing rows (years 1920-2008) with new GDP per Capita values
$ws.Range("E2").Value = "'1487"
$ws.Range("E3").Value = "'1466"
$ws.Range("E4").Value = "'1516"
$ws.Range("E5").Value = "'1554"
$ws.Range("E6").Value = "'1623"
$ws.Range("E7").Value = "'1476"
$ws.Range("E8").Value = "'1706"
$ws.Range("E9").Value = "'1478"
$ws.Range("E10").Value = "'1682"
$ws.Range("E11").Value = "'1661"
$ws.Range("E12").Value = "'1667"
$ws.Range("E13").Value = "'1474"
$ws.Range("E14").Value = "'1313"
$ws.Range("E15").Value = "'1471"
$ws.Range("E16").Value = "'1498"
$ws.Range("E17").Value = "'1629"
$ws.Range("E18").Value = "'1572"
$ws.Range("E19").Value = "'1699"
$ws.Range("E20").Value = "'1559"
$ws.Range("E21").Value = "'1651"
$ws.Range("E22").Value = "'1773"
$ws.Range("E23").Value = "'1714"
$ws.Range("E24").Value = "'1828"
$ws.Range("E25").Value = "'1970"
$ws.Range("E26").Value = "'1836"
$ws.Range("E27").Value = "'1741"
$ws.Range("E28").Value = "'1747"
$ws.Range("E29").Value = "'2173"
$ws.Range("E30").Value = "'2723"
$ws.Range("E31").Value = "'2432"
$ws.Range("E32").Value = "'2375"
$ws.Range("E33").Value = "'2361"
$ws.Range("E34").Value = "'2472"
$ws.Range("E35").Value = "'2577"
$ws.Range("E36").Value = "'2536"
$ws.Range("E37").Value = "'2590"
$ws.Range("E38").Value = "'2715"
$ws.Range("E39").Value = "'2775"
$ws.Range("E40").Value = "'2751"
$ws.Range("E41").Value = "'2788"
$ws.Range("E42").Value = "'2813"
$ws.Range("E43").Value = "'2821"
$ws.Range("E44").Value = "'3062"
$ws.Range("E45").Value = "'3096"
$ws.Range("E46").Value = "'3282"
$ws.Range("E47").Value = "'3351"
$ws.Range("E48").Value = "'3464"
$ws.Range("E49").Value = "'3534"
$ws.Range("E50").Value = "'3524"
$ws.Range("E51").Value = "'3520"
$ws.Range("E52").Value = "'3488"
$ws.Range("E53").Value = "'3553"
$ws.Range("E54").Value = "'3668"
$ws.Range("E55").Value = "'3767"
$ws.Range("E56").Value = "'3912"
$ws.Range("E57").Value = "'4030"
$ws.Range("E58").Value = "'4084"
$ws.Range("E59").Value = "'4221"
$ws.Range("E60").Value = "'4375"
$ws.Range("E61").Value = "'4194"
$ws.Range("E62").Value = "'3781"
$ws.Range("E63").Value = "'3507"
$ws.Range("E64").Value = "'3339"
$ws.Range("E65").Value = "'3330"
$ws.Range("E66").Value = "'3357"
$ws.Range("E67").Value = "'3366"
$ws.Range("E68").Value = "'3327"
$ws.Range("E69").Value = "'3351"
$ws.Range("E70").Value = "'3344"
$ws.Range("E71").Value = "'3322"
$ws.Range("E72").Value = "'3378"
$ws.Range("E73").Value = "'3493.57609938296"
$ws.Range("E74").Value = "'3750.51721814953"
$ws.Range("E75").Value = "'4052.95777781183"
$ws.Range("E76").Value = "'4316.39804713577"
$ws.Range("E77").Value = "'4592.28357203461"
$ws.Range("E78").Value = "'4671.33901624011"
$ws.Range("E79").Value = "'4882.68016685898"
$ws.Range("E80").Value = "'5077.8920578762"
$ws.Range("E81").Value = "'5263.88675536426"
$ws.Range("E82").Value = "'5403.7147358137"
$ws.Range("E83").Value = "'5550.18682788689"
$ws.Range("E84").Value = "'5751.06066526152"
$ws.Range("E85").Value = "'5957.20113760036"
$ws.Range("E86").Value = "'6145.97987183393"
$ws.Range("E87").Value = "'6441.49581866886"
$ws.Range("E88").Value = "'6782.23121360351"
$ws.Range("E89").Value = "'7139.18213958217"
$ws.Range("E90").Value = "'7312.65583506175"

# Add new rows for years 2009-2016
$ws.Range("A91").Value = 222
$ws.Range("B91").Value = "El Salvador"
$ws.Range("C91").Value = "GDP per Capita"
$ws.Range("D91").Value = 2009
$ws.Range("E91").Value = "'7166.15919423682"
$ws.Range("A92").Value = 222
$ws.Range("B92").Value = "El Salvador"
$ws.Range("C92").Value = "GDP per Capita"
$ws.Range("D92").Value = 2010
$ws.Range("E92").Value = "'7351.11455771672"
$ws.Range("A93").Value = 222
$ws.Range("B93").Value = "El Salvador"
$ws.Range("C93").Value = "GDP per Capita"
$ws.Range("D93").Value = 2011
$ws.Range("E93").Value = "'7607"
$ws.Range("A94").Value = 222
$ws.Range("B94").Value = "El Salvador"
$ws.Range("C94").Value = "GDP per Capita"
$ws.Range("D94").Value = 2012
$ws.Range("E94").Value = "'7726"
$ws.Range("A95").Value = 222
$ws.Range("B95").Value = "El Salvador"
$ws.Range("C95").Value = "GDP per Capita"
$ws.Range("D95").Value = 2013
$ws.Range("E95").Value = "'7845"
$ws.Range("A96").Value = 222
$ws.Range("B96").Value = "El Salvador"
$ws.Range("C96").Value = "GDP per Capita"
$ws.Range("D96").Value = 2014
$ws.Range("E96").Value = "'7935"
$ws.Range("A97").Value = 222
$ws.Range("B97").Value = "El Salvador"
$ws.Range("C97").Value = "GDP per Capita"
$ws.Range("D97").Value = 2015
$ws.Range("E97").Value = "'8109"
$ws.Range("A98").Value = 222
$ws.Range("B98").Value = "El Salvador"
$ws.Range("C98").Value = "GDP per Capita"
$ws.Range("D98").Value = 2016
$ws.Range("E98").Value = "'8280"
